# "updated import with regex"
# Fixes two text glitches introduced by the (regex-based) content importer:
#   - Slide 5: duplicated word "modificabili modificabili" -> "modificabili"
#   - Slide 8: typo "likeItellij" -> "likeIntellij"

$p = $ppt.ActivePresentation

function Fix-TextInPresentation {
    param(
        [string]$Search,
        [string]$Replacement
    )

    for ($si = 1; $si -le $p.Slides.Count; $si++) {
        $slide = $p.Slides.Item($si)
        for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
            $shape = $slide.Shapes.Item($shi)
            if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
                $tr = $shape.TextFrame.TextRange
                if ($tr.Text.Contains($Search)) {
                    $found = $tr.Find($Search)
                    while ($found -ne $null) {
                        $found.Text = $Replacement
                        $found = $tr.Find($Search)
                    }
                }
            }
        }
    }
}

Fix-TextInPresentation "Campi visibili o campi modificabili modificabili" "Campi visibili o campi modificabili"
Fix-TextInPresentation "likeItellij" "likeIntellij"
